# Generate Report for Handback
# Update the timestamp strings recorded during the handback status generation.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G2)
$wsOverview.Range("G2").Value = "2016-08-24 01:02:15"

# zh-cn sheet: "Correspond Handoff Datetime" (H2) and "Correspond Handback DateTime" (K2)
$wsZhCn.Range("H2").Value = "2016-08-24 01:02:10"
$wsZhCn.Range("K2").Value = "2016-08-24 01:02:35"

# de-de sheet: this datetime value is shared with the Overview sheet's G2 cell,
# so it must be kept in sync with the same new timestamp.
$wsDeDe.Range("H2").Value = "2016-08-24 01:02:15"
$wsDeDe.Range("K2").Value = "2016-08-24 01:02:42"
